# aggiornamento fino a 6/03
# Append three new daily rows (245-247) to the "nuovi positivi comuni" sheet,
# continuing the existing A1:AX244 table with the next three days.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each array is one new row: date serial in column A, then the per-comune
# counts through column AX (49 more values => 50 columns total, A..AX).
$row245 = @(44319,3,0,1,11,2,5,0,2,2,1,0,3,3,0,0,1,0,0,3,5,31,0,0,5,5,0,0,0,1,4,0,4,3,0,0,2,0,5,1,4,108,1,0,0,0,0,0,0,0)
$row246 = @(44320,0,1,0,4,4,1,1,0,0,0,1,1,0,0,0,0,1,0,0,2,14,0,2,0,1,2,0,0,1,1,1,5,2,1,1,0,0,0,0,3,51,0,0,1,0,0,0,0,0)
$row247 = @(44321,0,0,0,1,0,1,2,1,0,0,0,3,2,0,0,0,0,0,0,1,11,0,0,1,1,0,0,0,1,0,0,6,0,1,0,2,0,1,0,4,39,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $row245.Count; $i++) {
    $ws.Cells.Item(245, $i + 1).Value = $row245[$i]
}
for ($i = 0; $i -lt $row246.Count; $i++) {
    $ws.Cells.Item(246, $i + 1).Value = $row246[$i]
}
for ($i = 0; $i -lt $row247.Count; $i++) {
    $ws.Cells.Item(247, $i + 1).Value = $row247[$i]
}

# The date column (A) uses a dedicated cell style (border/bold/centered,
# custom date number format) throughout the table. Copy that formatting
# from the last existing data row (244) onto the three new date cells so
# they match the rest of the column instead of getting a brand-new style.
$ws.Cells.Item(244, 1).Copy() | Out-Null
$ws.Cells.Item(245, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(244, 1).Copy() | Out-Null
$ws.Cells.Item(246, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(244, 1).Copy() | Out-Null
$ws.Cells.Item(247, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
